$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("April")

# Update row labels in column A for rows 39-43 (library name reordering)
$ws.Cells.Item(39, 1).Value = "Prairie Hills Schools - Axtell Public School"
$ws.Cells.Item(40, 1).Value = "Prairie Hills Schools - Sabetha Elementary School"
$ws.Cells.Item(41, 1).Value = "Prairie Hills Schools - Sabetha High School"
$ws.Cells.Item(42, 1).Value = "Prairie Hills Schools - Sabetha Middle School"
$ws.Cells.Item(43, 1).Value = "Pomona Community Library"

# Populate borrow/lend statistics for each library (rows 2-54)
$ws.Cells.Item(2, 2).Value = 1390
$ws.Cells.Item(2, 3).Value = 1244
$ws.Cells.Item(2, 4).Value = 146
$ws.Cells.Item(2, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(2, 7).Value = "1.12 : 1"

$ws.Cells.Item(3, 2).Value = 524
$ws.Cells.Item(3, 3).Value = 456
$ws.Cells.Item(3, 4).Value = 68
$ws.Cells.Item(3, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(3, 7).Value = "1.15 : 1"

$ws.Cells.Item(4, 2).Value = 1129
$ws.Cells.Item(4, 3).Value = 1033
$ws.Cells.Item(4, 4).Value = 96
$ws.Cells.Item(4, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(4, 7).Value = "1.09 : 1"

$ws.Cells.Item(5, 2).Value = 30
$ws.Cells.Item(5, 3).Value = 161
$ws.Cells.Item(5, 4).Value = -131
$ws.Cells.Item(5, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(5, 7).Value = "0.19 : 1"

$ws.Cells.Item(6, 2).Value = 1132
$ws.Cells.Item(6, 3).Value = 1219
$ws.Cells.Item(6, 4).Value = -87
$ws.Cells.Item(6, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(6, 7).Value = "0.93 : 1"

$ws.Cells.Item(7, 2).Value = 211
$ws.Cells.Item(7, 3).Value = 134
$ws.Cells.Item(7, 4).Value = 77
$ws.Cells.Item(7, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(7, 7).Value = "1.57 : 1"

$ws.Cells.Item(8, 2).Value = 117
$ws.Cells.Item(8, 3).Value = 168
$ws.Cells.Item(8, 4).Value = -51
$ws.Cells.Item(8, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(8, 7).Value = "0.70 : 1"

$ws.Cells.Item(9, 2).Value = 46
$ws.Cells.Item(9, 3).Value = 70
$ws.Cells.Item(9, 4).Value = -24
$ws.Cells.Item(9, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(9, 7).Value = "0.66 : 1"

$ws.Cells.Item(10, 2).Value = 0
$ws.Cells.Item(10, 3).Value = 46
$ws.Cells.Item(10, 4).Value = -46
$ws.Cells.Item(10, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(10, 7).Value = "0.00 : 1"

$ws.Cells.Item(11, 2).Value = 0
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 0

$ws.Cells.Item(12, 2).Value = 20
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 10
$ws.Cells.Item(12, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(12, 7).Value = "2.00 : 1"

$ws.Cells.Item(13, 2).Value = 77
$ws.Cells.Item(13, 3).Value = 57
$ws.Cells.Item(13, 4).Value = 20
$ws.Cells.Item(13, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(13, 7).Value = "1.35 : 1"

$ws.Cells.Item(14, 2).Value = 254
$ws.Cells.Item(14, 3).Value = 220
$ws.Cells.Item(14, 4).Value = 34
$ws.Cells.Item(14, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(14, 7).Value = "1.15 : 1"

$ws.Cells.Item(15, 2).Value = 56
$ws.Cells.Item(15, 3).Value = 121
$ws.Cells.Item(15, 4).Value = -65
$ws.Cells.Item(15, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(15, 7).Value = "0.46 : 1"

$ws.Cells.Item(16, 2).Value = 66
$ws.Cells.Item(16, 3).Value = 140
$ws.Cells.Item(16, 4).Value = -74
$ws.Cells.Item(16, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(16, 7).Value = "0.47 : 1"

$ws.Cells.Item(17, 2).Value = 620
$ws.Cells.Item(17, 3).Value = 469
$ws.Cells.Item(17, 4).Value = 151
$ws.Cells.Item(17, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(17, 7).Value = "1.32 : 1"

$ws.Cells.Item(18, 2).Value = 71
$ws.Cells.Item(18, 3).Value = 96
$ws.Cells.Item(18, 4).Value = -25
$ws.Cells.Item(18, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(18, 7).Value = "0.74 : 1"

$ws.Cells.Item(19, 2).Value = 557
$ws.Cells.Item(19, 3).Value = 472
$ws.Cells.Item(19, 4).Value = 85
$ws.Cells.Item(19, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(19, 7).Value = "1.18 : 1"

$ws.Cells.Item(20, 2).Value = 33
$ws.Cells.Item(20, 3).Value = 67
$ws.Cells.Item(20, 4).Value = -34
$ws.Cells.Item(20, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(20, 7).Value = "0.49 : 1"

$ws.Cells.Item(21, 2).Value = 507
$ws.Cells.Item(21, 3).Value = 342
$ws.Cells.Item(21, 4).Value = 165
$ws.Cells.Item(21, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(21, 7).Value = "1.48 : 1"

$ws.Cells.Item(22, 2).Value = 90
$ws.Cells.Item(22, 3).Value = 233
$ws.Cells.Item(22, 4).Value = -143
$ws.Cells.Item(22, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(22, 7).Value = "0.39 : 1"

$ws.Cells.Item(23, 2).Value = 616
$ws.Cells.Item(23, 3).Value = 490
$ws.Cells.Item(23, 4).Value = 126
$ws.Cells.Item(23, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(23, 7).Value = "1.26 : 1"

$ws.Cells.Item(24, 2).Value = 1710
$ws.Cells.Item(24, 3).Value = 1144
$ws.Cells.Item(24, 4).Value = 566
$ws.Cells.Item(24, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(24, 7).Value = "1.49 : 1"

$ws.Cells.Item(25, 2).Value = 203
$ws.Cells.Item(25, 3).Value = 441
$ws.Cells.Item(25, 4).Value = -238
$ws.Cells.Item(25, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(25, 7).Value = "0.46 : 1"

$ws.Cells.Item(26, 2).Value = 0
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(26, 4).Value = 0

$ws.Cells.Item(27, 2).Value = 201
$ws.Cells.Item(27, 3).Value = 168
$ws.Cells.Item(27, 4).Value = 33
$ws.Cells.Item(27, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(27, 7).Value = "1.20 : 1"

$ws.Cells.Item(28, 2).Value = 94
$ws.Cells.Item(28, 3).Value = 92
$ws.Cells.Item(28, 4).Value = 2
$ws.Cells.Item(28, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(28, 7).Value = "1.02 : 1"

$ws.Cells.Item(29, 2).Value = 571
$ws.Cells.Item(29, 3).Value = 425
$ws.Cells.Item(29, 4).Value = 146
$ws.Cells.Item(29, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(29, 7).Value = "1.34 : 1"

$ws.Cells.Item(30, 2).Value = 38
$ws.Cells.Item(30, 3).Value = 24
$ws.Cells.Item(30, 4).Value = 14
$ws.Cells.Item(30, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(30, 7).Value = "1.58 : 1"

$ws.Cells.Item(31, 2).Value = 78
$ws.Cells.Item(31, 3).Value = 253
$ws.Cells.Item(31, 4).Value = -175
$ws.Cells.Item(31, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(31, 7).Value = "0.31 : 1"

$ws.Cells.Item(32, 2).Value = 371
$ws.Cells.Item(32, 3).Value = 580
$ws.Cells.Item(32, 4).Value = -209
$ws.Cells.Item(32, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(32, 7).Value = "0.64 : 1"

$ws.Cells.Item(33, 2).Value = 443
$ws.Cells.Item(33, 3).Value = 417
$ws.Cells.Item(33, 4).Value = 26
$ws.Cells.Item(33, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(33, 7).Value = "1.06 : 1"

$ws.Cells.Item(34, 2).Value = 193
$ws.Cells.Item(34, 3).Value = 149
$ws.Cells.Item(34, 4).Value = 44
$ws.Cells.Item(34, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(34, 7).Value = "1.30 : 1"

$ws.Cells.Item(35, 2).Value = 787
$ws.Cells.Item(35, 3).Value = 1256
$ws.Cells.Item(35, 4).Value = -469
$ws.Cells.Item(35, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(35, 7).Value = "0.63 : 1"

$ws.Cells.Item(36, 2).Value = 172
$ws.Cells.Item(36, 3).Value = 541
$ws.Cells.Item(36, 4).Value = -369
$ws.Cells.Item(36, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(36, 7).Value = "0.32 : 1"

$ws.Cells.Item(37, 2).Value = 485
$ws.Cells.Item(37, 3).Value = 271
$ws.Cells.Item(37, 4).Value = 214
$ws.Cells.Item(37, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(37, 7).Value = "1.79 : 1"

$ws.Cells.Item(38, 2).Value = 16
$ws.Cells.Item(38, 3).Value = 151
$ws.Cells.Item(38, 4).Value = -135
$ws.Cells.Item(38, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(38, 7).Value = "0.11 : 1"

$ws.Cells.Item(39, 2).Value = 15
$ws.Cells.Item(39, 3).Value = 100
$ws.Cells.Item(39, 4).Value = -85
$ws.Cells.Item(39, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(39, 7).Value = "0.15 : 1"

$ws.Cells.Item(40, 2).Value = 72
$ws.Cells.Item(40, 3).Value = 89
$ws.Cells.Item(40, 4).Value = -17
$ws.Cells.Item(40, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(40, 7).Value = "0.81 : 1"

$ws.Cells.Item(41, 2).Value = 10
$ws.Cells.Item(41, 3).Value = 22
$ws.Cells.Item(41, 4).Value = -12
$ws.Cells.Item(41, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(41, 7).Value = "0.45 : 1"

$ws.Cells.Item(42, 2).Value = 8
$ws.Cells.Item(42, 3).Value = 62
$ws.Cells.Item(42, 4).Value = -54
$ws.Cells.Item(42, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(42, 7).Value = "0.13 : 1"

$ws.Cells.Item(43, 2).Value = 60
$ws.Cells.Item(43, 3).Value = 86
$ws.Cells.Item(43, 4).Value = -26
$ws.Cells.Item(43, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(43, 7).Value = "0.70 : 1"

$ws.Cells.Item(44, 2).Value = 75
$ws.Cells.Item(44, 3).Value = 229
$ws.Cells.Item(44, 4).Value = -154
$ws.Cells.Item(44, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(44, 7).Value = "0.33 : 1"

$ws.Cells.Item(45, 2).Value = 606
$ws.Cells.Item(45, 3).Value = 541
$ws.Cells.Item(45, 4).Value = 65
$ws.Cells.Item(45, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(45, 7).Value = "1.12 : 1"

$ws.Cells.Item(46, 2).Value = 1122
$ws.Cells.Item(46, 3).Value = 645
$ws.Cells.Item(46, 4).Value = 477
$ws.Cells.Item(46, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(46, 7).Value = "1.74 : 1"

$ws.Cells.Item(47, 2).Value = 233
$ws.Cells.Item(47, 3).Value = 570
$ws.Cells.Item(47, 4).Value = -337
$ws.Cells.Item(47, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(47, 7).Value = "0.41 : 1"

$ws.Cells.Item(48, 2).Value = 658
$ws.Cells.Item(48, 3).Value = 231
$ws.Cells.Item(48, 4).Value = 427
$ws.Cells.Item(48, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(48, 7).Value = "2.85 : 1"

$ws.Cells.Item(49, 2).Value = 614
$ws.Cells.Item(49, 3).Value = 542
$ws.Cells.Item(49, 4).Value = 72
$ws.Cells.Item(49, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(49, 7).Value = "1.13 : 1"

$ws.Cells.Item(50, 2).Value = 163
$ws.Cells.Item(50, 3).Value = 169
$ws.Cells.Item(50, 4).Value = -6
$ws.Cells.Item(50, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(50, 7).Value = "0.96 : 1"

$ws.Cells.Item(51, 2).Value = 371
$ws.Cells.Item(51, 3).Value = 388
$ws.Cells.Item(51, 4).Value = -17
$ws.Cells.Item(51, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(51, 7).Value = "0.96 : 1"

$ws.Cells.Item(52, 2).Value = 163
$ws.Cells.Item(52, 3).Value = 184
$ws.Cells.Item(52, 4).Value = -21
$ws.Cells.Item(52, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(52, 7).Value = "0.89 : 1"

$ws.Cells.Item(53, 2).Value = 51
$ws.Cells.Item(53, 3).Value = 255
$ws.Cells.Item(53, 4).Value = -204
$ws.Cells.Item(53, 6).Value = "We lent more than we borrowed"
$ws.Cells.Item(53, 7).Value = "0.20 : 1"

$ws.Cells.Item(54, 2).Value = 365
$ws.Cells.Item(54, 3).Value = 221
$ws.Cells.Item(54, 4).Value = 144
$ws.Cells.Item(54, 5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(54, 7).Value = "1.65 : 1"
